$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 32581
$ws.Range("K2").Value = 32581

$ws.Range("I3").Value = 11.01169489224504
$ws.Range("J3").Value = 0.2181639605905282
$ws.Range("K3").Value = 0.1702034928332464

$ws.Range("I4").Value = 3.240459464955947
$ws.Range("J4").Value = 0.4130056685601487
$ws.Range("K4").Value = 0.1067817563449238

$ws.Range("I5").Value = 5.42
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0

$ws.Range("I6").Value = 7.9
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.09

$ws.Range("I7").Value = 10.99
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.15

$ws.Range("I8").Value = 13.47
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.23

$ws.Range("I9").Value = 23.22
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 0.83

$ws.Range("I10").Value = 0.2085503016908977
$ws.Range("J10").Value = 1.364888487327168
$ws.Range("K10").Value = 1.064668636768324

$ws.Range("I11").Value = -0.6716091079813706
$ws.Range("J11").Value = -0.137087836131323
$ws.Range("K11").Value = 1.223686678285682

$ws.Range("I12").Value = 6.03
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0.04

$ws.Range("I13").Value = 16.32
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 0.38

$ws.Range("K14").Value = 0

$ws.Range("K15").Value = 0

$ws.Range("K16").Value = 77

$ws.Range("J17").Value = 0.006138547005923698
$ws.Range("K17").Value = 0.2363340597280624

$ws.Range("J18").Value = 25473
$ws.Range("K18").Value = 9

$ws.Range("J19").Value = 78.18360394094718
$ws.Range("K19").Value = 0.02762346152665664

$ws.Range("I20").Value = 17.8
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 0.83

$ws.Range("I21").Value = 5.57
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0.14

$ws.Range("I22").Value = 10.50057754402258
$ws.Range("J22").Value = 0.1705736822628154
$ws.Range("K22").Value = 0.01140234348810668

$ws.Range("I23").Value = 324459.59
$ws.Range("J23").Value = 7108
$ws.Range("K23").Value = 5545.400000000001
